$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header text updates (rich-text shared strings) ----
$ws.Range("A8").Value = "Volume 29   Number  40"
$ws.Range("C9").Value = "Report Covering the Week  10/3/2022  Through  10/9/2022"

# ---- Cells changing from numeric style to text-placeholder style (target style 14) ----
$textCells = [ordered]@{
    F15 = "0"
    C17 = "0"
    C18 = "0"
    C22 = "0"
    D22 = "0"
    E22 = "***.*"
    C26 = "0"
    G30 = "0"
    H30 = "***.*"
}
foreach ($ref in $textCells.Keys) {
    $ws.Range($ref).NumberFormat = "@"
    $ws.Range($ref).Value = $textCells[$ref]
    $ws.Range("A14").Copy()
    $ws.Range($ref).PasteSpecial(-4122)
}

# ---- Cells changing from text-placeholder style to a numeric style (target style 15 or 16) ----
$numCells15 = [ordered]@{
    D15 = 1
    G15 = 1
    D17 = 2
    D18 = 1
    D23 = 1
    G23 = 1
    J23 = 1
    C27 = 1
}
foreach ($ref in $numCells15.Keys) {
    $ws.Range("I14").Copy()
    $ws.Range($ref).PasteSpecial(-4122)
    $ws.Range($ref).Value = $numCells15[$ref]
}

$numCells16 = [ordered]@{
    E15 = -100
    H15 = -100
    E17 = -100
    E18 = -100
    E23 = -100
    H23 = -100
    K23 = 300
}
foreach ($ref in $numCells16.Keys) {
    $ws.Range("M14").Copy()
    $ws.Range($ref).PasteSpecial(-4122)
    $ws.Range($ref).Value = $numCells16[$ref]
}

# ---- Plain value updates (style unchanged) ----
$plainCells = [ordered]@{
    J15 = 8
    K15 = 12.5
    L15 = 50
    M15 = 200
    C16 = 3
    F16 = 6
    G16 = 4
    H16 = 50
    I16 = 69
    J16 = 33
    K16 = 109.090909090909
    L16 = 50
    M16 = -23.333333333333
    N16 = -86.116700201207
    F17 = 2
    G17 = 3
    H17 = -33.333333333333
    J17 = 47
    K17 = 38.297872340425
    L17 = 71.052631578947
    M17 = 41.304347826087
    N17 = -34.343434343434
    F18 = 13
    H18 = 116.666666666667
    I18 = 86
    J18 = 69
    K18 = 24.637681159420
    L18 = -4.444444444444
    M18 = -13.131313131313
    N18 = -91.871455576559
    C19 = 10
    D19 = 10
    E19 = 0
    F19 = 38
    G19 = 38
    H19 = 0
    I19 = 397
    J19 = 270
    K19 = 47.037037037037
    L19 = 76.444444444444
    M19 = 32.333333333333
    N19 = -49.102564102564
    D20 = 1
    E20 = 100
    F20 = 8
    G20 = 10
    H20 = -20
    I20 = 78
    J20 = 53
    K20 = 47.169811320754
    L20 = 56
    M20 = 0
    N20 = -97.077557137504
    C21 = 15
    D21 = 15
    E21 = 0
    F21 = 67
    G21 = 62
    H21 = 8.064516129032
    I21 = 706
    J21 = 480
    K21 = 47.083333333333
    L21 = 55.164835164835
    M21 = 14.424635332252
    N21 = -86.19475948377
    G22 = 2
    H22 = 50
    J22 = 10
    K22 = 180
    L22 = 100
    M22 = 55.555555555555
    C24 = 37
    D24 = 26
    E24 = 42.307692307692
    F24 = 130
    G24 = 116
    H24 = 12.068965517241
    I24 = 1395
    J24 = 977
    K24 = 42.784032753326
    L24 = 76.136363636363
    M24 = 84.280052840158
    C25 = 3
    D25 = 1
    E25 = 200
    F25 = 9
    G25 = 10
    H25 = -10
    I25 = 149
    J25 = 142
    K25 = 4.929577464788
    L25 = 50.505050505050
    M25 = -6.875
    E26 = -100
    F26 = 1
    H26 = -50
    J26 = 14
    K26 = 0
    L26 = -26.315789473684
    F27 = 4
    G27 = 2
    H27 = 100
    I27 = 34
    J27 = 25
    K27 = 36
    L27 = -10.526315789473
    N28 = -66.666666666666
    N29 = -57.142857142857
}
foreach ($ref in $plainCells.Keys) {
    $ws.Range($ref).Value = $plainCells[$ref]
}
